$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# Row 54: new "7ème semaine" section header, formatted like the
# previous week-header rows (A41, A48) and merged across A:C.
# ------------------------------------------------------------------
$ws.Range("A48:C48").Copy()
$ws.Range("A54:C54").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false
$ws.Range("A54:C54").Merge()
$ws.Cells.Item(54,1).Value = "7ème semaine "

# ------------------------------------------------------------------
# Row 55: new activity entry (copy formatting from row 49, which has
# the same style pattern: date / wrapped text / duration)
# ------------------------------------------------------------------
$ws.Range("A49:C49").Copy()
$ws.Range("A55:C55").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Cells.Item(55,1).Value = 43186
$ws.Cells.Item(55,2).Value = "J'ai eu un problème quand j'ajoutais un article dans un panier d'un utilisateur fraîchement créé. Mon article ne s'ajoutait pas la première fois que j'appuyais sur le bouton mais la deuxième fois. J'avais ajouté un un champ en trop dans ma requête qui faisait que ça ne fonctionnait pas. "
$ws.Cells.Item(55,3).Value = "1 période"
$ws.Rows.Item(55).RowHeight = 60

# ------------------------------------------------------------------
# Row 56
# ------------------------------------------------------------------
$ws.Range("A49:C49").Copy()
$ws.Range("A56:C56").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Cells.Item(56,1).Value = 43186
$ws.Cells.Item(56,2).Value = "J'ai aussi ajouté des contrôles lors de la connexion et de l'inscription, lorsque le mot de passe et faux à la connexion un message apparait. A l'inscription si le login est déjà utilisé par quelqu'un d'autre un message apparait."
$ws.Cells.Item(56,3).Value = "1 période"
$ws.Rows.Item(56).RowHeight = 45

# ------------------------------------------------------------------
# Row 57
# ------------------------------------------------------------------
$ws.Range("A49:C49").Copy()
$ws.Range("A57:C57").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Cells.Item(57,1).Value = 43186
$ws.Cells.Item(57,2).Value = "M. Carrel est venu m'aider parce que j'avais un problème avec mon panier. Lorsque j'ajoutais et supprimait un article dans le panier d'un utilisateur, je devais rafraichir à nouveau la page pour que dans l'icône en haut à droite de ma page s'affiche mes articles dans le panier et le nombre que j'ai d'article dedans. J'ai changé de place un require_once pour que mon affichage de mon menu se fassent après le traitement des données. "
$ws.Cells.Item(57,3).Value = "1 période"
$ws.Rows.Item(57).RowHeight = 75

# ------------------------------------------------------------------
# Scroll/selection state, like in the edited workbook
# ------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 52
$ws.Range("C58").Select()
